$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 13:05"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 1385893
$ws.Range("C4").Value = 59
$ws.Range("E4").Value = 1041872
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 81796

# Row 13: Iran - refreshed totals
$ws.Range("B13").Value = 110767
$ws.Range("C13").Value = 1481
$ws.Range("D13").Value = 88357
$ws.Range("E13").Value = 15677
$ws.Range("F13").Value = 2713
$ws.Range("G13").Value = 48
$ws.Range("H13").Value = 6733

# Row 28: Singapur - refreshed totals
$ws.Range("C28").Value = 884

# Row 55: Marruecos - refreshed totals
$ws.Range("B55").Value = 6380
$ws.Range("C55").Value = 99
$ws.Range("D55").Value = 2930
$ws.Range("E55").Value = 3262

# Row 57: Finlandia - refreshed totals
$ws.Range("B57").Value = 6003
$ws.Range("C57").Value = 19
$ws.Range("E57").Value = 1732

# Row 60: Barein - refreshed totals
$ws.Range("E60").Value = 3075
$ws.Range("F60").Value = 5
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 9

# Rows 62-63: Afganistan overtakes Ghana in total cases, so they swap order
$ws.Range("A62").Value = "Afganistan"
$ws.Range("B62").Value = 4963
$ws.Range("C62").Value = 276
$ws.Range("D62").Value = 610
$ws.Range("E62").Value = 4226
$ws.Range("F62").Value = 7
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 127

$ws.Range("A63").Value = "Ghana"
$ws.Range("B63").Value = 4700
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 494
$ws.Range("E63").Value = 4184
$ws.Range("F63").Value = 5
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 22

# Row 87: Sudan - refreshed totals
$ws.Range("D87").Value = 161
$ws.Range("E87").Value = 1291

# Row 104: Albania - refreshed totals
$ws.Range("F104").Value = 1

# Rows 215-216: San Bartolome and Sahara Occidental tied, order swaps
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"
